# Added Test Cases and data for VLink Card Not Accepted Error.
# Updates the "Date" column (column B) for the existing 21 data rows
# with new timestamp values, as new test case runs were appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDates = @(
    "Wed Feb 08 16:46:36 EST 2023",
    "Wed Feb 08 16:47:25 EST 2023",
    "Wed Feb 08 16:48:19 EST 2023",
    "Wed Feb 08 16:49:06 EST 2023",
    "Wed Feb 08 16:49:51 EST 2023",
    "Wed Feb 08 16:50:37 EST 2023",
    "Wed Feb 08 16:51:23 EST 2023",
    "Wed Feb 08 16:52:16 EST 2023",
    "Wed Feb 08 16:53:08 EST 2023",
    "Wed Feb 08 16:54:00 EST 2023",
    "Wed Feb 08 16:54:49 EST 2023",
    "Wed Feb 08 16:55:43 EST 2023",
    "Wed Feb 08 16:56:29 EST 2023",
    "Wed Feb 08 16:57:20 EST 2023",
    "Wed Feb 08 16:58:10 EST 2023",
    "Wed Feb 08 16:59:00 EST 2023",
    "Wed Feb 08 16:59:52 EST 2023",
    "Wed Feb 08 17:00:43 EST 2023",
    "Wed Feb 08 17:01:37 EST 2023",
    "Wed Feb 08 17:02:30 EST 2023",
    "Wed Feb 08 17:03:23 EST 2023"
)

$startRow = 2
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
